$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new time log entry: 09/15/2023 - Internship - Completed 8 hours assisting with daily operations
# Copy formatting from the row above first so the new row matches existing style/number format
$ws.Range("A4:C4").Copy()
$ws.Range("A5:C5").PasteSpecial(-4122)

# Set the date as a raw serial number (avoids Excel auto-applying a new date format)
$ws.Range("A5").Value = 45184
$ws.Range("B5").Value = $ws.Range("B4").Value2
$ws.Range("C5").Value = $ws.Range("C4").Value2

$ws.Range("C9").Select()
